$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain value updates (non-numeric-looking strings; Excel keeps these as text automatically)
$plainItems = @(
    @{Ref="D2"; Text='67.127.84'},
    @{Ref="E2"; Text='  -0.63%  '},
    @{Ref="D3"; Text='3.108.36'},
    @{Ref="E3"; Text='  -0.39%  '},
    @{Ref="E4"; Text='  -0.06%  '},
    @{Ref="E5"; Text='  -0.92%  '},
    @{Ref="E6"; Text='  +2.43%  '},
    @{Ref="E7"; Text='  -0.04%  '},
    @{Ref="D8"; Text='3.106.32'},
    @{Ref="E8"; Text='  -0.32%  '},
    @{Ref="E9"; Text='  -1.31%  '},
    @{Ref="E10"; Text='  -1.93%  '},
    @{Ref="E11"; Text='  -1.02%  '},
    @{Ref="E12"; Text='  -2.06%  '},
    @{Ref="E13"; Text='  -3.02%  '},
    @{Ref="E14"; Text='  -2.19%  '},
    @{Ref="E15"; Text='  -0.09%  '},
    @{Ref="D16"; Text='3.622.14'},
    @{Ref="E16"; Text='  -0.46%  '},
    @{Ref="D17"; Text='67.062.01'},
    @{Ref="E17"; Text='  -0.67%  '},
    @{Ref="E18"; Text='  -1.03%  '},
    @{Ref="D19"; Text='3.103.65'},
    @{Ref="E19"; Text='  -0.54%  '},
    @{Ref="E21"; Text='  +0.68%  '},
    @{Ref="E22"; Text='  -0.58%  '},
    @{Ref="E23"; Text='  -2.05%  '},
    @{Ref="E24"; Text='  -0.37%  '},
    @{Ref="E25"; Text='  -3.89%  '},
    @{Ref="E26"; Text='  -1.77%  '},
    @{Ref="E27"; Text='  -4.30%  '},
    @{Ref="E28"; Text='  +0.03%  '},
    @{Ref="E29"; Text='  +0.18%  '},
    @{Ref="E30"; Text='  -1.63%  '},
    @{Ref="E31"; Text='  -2.86%  '},
    @{Ref="E32"; Text='  -0.79%  '},
    @{Ref="E33"; Text='  -1.84%  '},
    @{Ref="D34"; Text='0.0₃0942'},
    @{Ref="E34"; Text='  -0.50%  '},
    @{Ref="E35"; Text='  -0.09%  '},
    @{Ref="E36"; Text='  +1.19%  '},
    @{Ref="B37"; Text='Filecoin'},
    @{Ref="C37"; Text='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'},
    @{Ref="E37"; Text='  -3.95%  '},
    @{Ref="B38"; Text='Mantle'},
    @{Ref="C38"; Text='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'},
    @{Ref="E38"; Text='  -2.63%  '},
    @{Ref="E39"; Text='  +1.21%  '},
    @{Ref="E40"; Text='  -1.85%  '},
    @{Ref="E41"; Text='  -1.19%  '},
    @{Ref="E42"; Text='  -0.93%  '},
    @{Ref="E43"; Text='  +5.56%  '},
    @{Ref="E44"; Text='  -2.08%  '},
    @{Ref="D45"; Text='2.791.58'},
    @{Ref="E45"; Text='  -0.47%  '},
    @{Ref="E46"; Text='  -3.63%  '},
    @{Ref="E47"; Text='  -1.98%  '},
    @{Ref="E48"; Text='  -0.20%  '},
    @{Ref="E49"; Text='  +0.02%  '},
    @{Ref="E50"; Text='  +1.64%  '},
    @{Ref="E51"; Text='  +3.90%  '}
)
foreach ($item in $plainItems) {
    $ws.Range($item.Ref).Value = $item.Text
}

# Numeric-looking strings must be forced to remain text (matches original inlineStr cells),
# otherwise Excel auto-converts them to real numbers. Temporarily set a Text number format,
# assign the value, then restore the original style so the cell style is unchanged.
$textItems = @(
    @{Ref="D5"; Text='575.10'},
    @{Ref="D6"; Text='178.22'},
    @{Ref="D10"; Text='6.38'},
    @{Ref="D12"; Text='0.469'},
    @{Ref="D13"; Text='0.0000241'},
    @{Ref="D14"; Text='36.09'},
    @{Ref="D18"; Text='7.05'},
    @{Ref="D20"; Text='16.72'},
    @{Ref="D21"; Text='496.74'},
    @{Ref="D22"; Text='7.78'},
    @{Ref="D23"; Text='0.688'},
    @{Ref="D24"; Text='83.73'},
    @{Ref="D25"; Text='12.62'},
    @{Ref="D27"; Text='10.10'},
    @{Ref="D28"; Text='1.00'},
    @{Ref="D29"; Text='7.92'},
    @{Ref="D30"; Text='2.31'},
    @{Ref="D32"; Text='28.22'},
    @{Ref="D35"; Text='0.999'},
    @{Ref="D36"; Text='47.49'},
    @{Ref="D37"; Text='5.60'},
    @{Ref="D38"; Text='0.945'},
    @{Ref="D40"; Text='49.09'},
    @{Ref="D41"; Text='2.02'},
    @{Ref="D43"; Text='2.74'},
    @{Ref="D44"; Text='8.31'},
    @{Ref="D46"; Text='372.25'},
    @{Ref="D47"; Text='0.0345'},
    @{Ref="D48"; Text='135.37'},
    @{Ref="D50"; Text='25.48'},
    @{Ref="D51"; Text='2.29'}
)
foreach ($item in $textItems) {
    $cell = $ws.Range($item.Ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $item.Text
    $cell.Style = $origStyle
}
